# Updates cryptos list values (price/volume columns, and the row shift
# caused by inserting a new 'BitDAO' entry at row 24) per the commit diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'30.434.09"
$ws.Range("E2").Value = "'  -0.42%  "

$ws.Range("D3").Value = "'1.862.26"
$ws.Range("E3").Value = "'  -0.80%  "

$ws.Range("D4").Value = "'0.9997"
$ws.Range("E4").Value = "'  -0.10%  "

$ws.Range("D5").Value = "'234.91"
$ws.Range("E5").Value = "'  -1.48%  "

$ws.Range("D6").Value = "'0.9992"
$ws.Range("E6").Value = "'  -0.11%  "

$ws.Range("D7").Value = "'0.4743"
$ws.Range("E7").Value = "'  -1.11%  "

$ws.Range("D8").Value = "'0.2752"
$ws.Range("E8").Value = "'  -2.15%  "

$ws.Range("D9").Value = "'0.06445"
$ws.Range("E9").Value = "'  -0.82%  "

$ws.Range("D10").Value = "'1.879.07"
$ws.Range("E10").Value = "'  -1.02%  "

$ws.Range("D11").Value = "'0.07449"
$ws.Range("E11").Value = "'  -0.42%  "

$ws.Range("D12").Value = "'16.37"
$ws.Range("E12").Value = "'  -0.73%  "

$ws.Range("D13").Value = "'5.002"
$ws.Range("E13").Value = "'  -1.51%  "

$ws.Range("D14").Value = "'85.84"
$ws.Range("E14").Value = "'  -2.22%  "

$ws.Range("D15").Value = "'0.6358"
$ws.Range("E15").Value = "'  -3.93%  "

$ws.Range("D16").Value = "'30.390.57"
$ws.Range("E16").Value = "'  -0.44%  "

$ws.Range("E17").Value = "'  -0.04%  "

$ws.Range("D18").Value = "'231.20"
$ws.Range("E18").Value = "'  +1.93%  "

$ws.Range("E19").Value = "'  -3.13%  "

$ws.Range("E20").Value = "'  -1.73%  "

$ws.Range("D21").Value = "'2.103.04"
$ws.Range("E21").Value = "'  -2.65%  "

$ws.Range("D22").Value = "'0.9997"
$ws.Range("E22").Value = "'  -0.18%  "

$ws.Range("D23").Value = "'5.024"
$ws.Range("E23").Value = "'  -4.53%  "

$ws.Range("B24").Value = "'BitDAO"
$ws.Range("C24").Value = "'https://coinranking.com/coin/N2IgQ9Xme+bitdao-bit"
$ws.Range("D24").Value = "'0.3908"
$ws.Range("E24").Value = "'  +3.35%  "

$ws.Range("B25").Value = "'Chainlink"
$ws.Range("C25").Value = "'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D25").Value = "'6.016"
$ws.Range("E25").Value = "'  -2.19%  "

$ws.Range("B26").Value = "'Cosmos"
$ws.Range("C26").Value = "'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D26").Value = "'9.264"
$ws.Range("E26").Value = "'  -0.27%  "

$ws.Range("B27").Value = "'Monero"
$ws.Range("C27").Value = "'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D27").Value = "'166.39"
$ws.Range("E27").Value = "'  -0.39%  "

$ws.Range("B28").Value = "'EthereumClassic"
$ws.Range("C28").Value = "'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D28").Value = "'17.98"
$ws.Range("E28").Value = "'  -2.26%  "

$ws.Range("B29").Value = "'LidoDAOToken"
$ws.Range("C29").Value = "'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D29").Value = "'1.900"
$ws.Range("E29").Value = "'  -1.86%  "

$ws.Range("B30").Value = "'Stellar"
$ws.Range("C30").Value = "'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D30").Value = "'0.1043"
$ws.Range("E30").Value = "'  +7.52%  "

$ws.Range("B31").Value = "'Toncoin"
$ws.Range("C31").Value = "'https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D31").Value = "'1.399"
$ws.Range("E31").Value = "'  +0.04%  "

$ws.Range("B32").Value = "'InternetComputer(DFINITY)"
$ws.Range("C32").Value = "'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D32").Value = "'4.167"
$ws.Range("E32").Value = "'  -3.72%  "

$ws.Range("B33").Value = "'Filecoin"
$ws.Range("C33").Value = "'https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D33").Value = "'3.946"
$ws.Range("E33").Value = "'  -1.22%  "

$ws.Range("B34").Value = "'Hedera"
$ws.Range("C34").Value = "'https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D34").Value = "'0.04932"
$ws.Range("E34").Value = "'  -2.19%  "

$ws.Range("B35").Value = "'ARBITRUM"
$ws.Range("C35").Value = "'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D35").Value = "'1.162"
$ws.Range("E35").Value = "'  -4.67%  "

$ws.Range("B36").Value = "'ImmutableX"
$ws.Range("C36").Value = "'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D36").Value = "'0.7273"
$ws.Range("E36").Value = "'  -2.57%  "

$ws.Range("B37").Value = "'Frax"
$ws.Range("C37").Value = "'https://coinranking.com/coin/KfWtaeV1W+frax-frax"
$ws.Range("D37").Value = "'0.9995"
$ws.Range("E37").Value = "'  -0.15%  "

$ws.Range("B38").Value = "'HuobiToken"
$ws.Range("C38").Value = "'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D38").Value = "'2.697"
$ws.Range("E38").Value = "'  -0.45%  "

$ws.Range("B39").Value = "'VeChain"
$ws.Range("C39").Value = "'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D39").Value = "'0.01886"
$ws.Range("E39").Value = "'  +1.75%  "

$ws.Range("B40").Value = "'MXToken"
$ws.Range("C40").Value = "'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D40").Value = "'2.651"
$ws.Range("E40").Value = "'  +0.45%  "

$ws.Range("B41").Value = "'TrustWalletToken"
$ws.Range("C41").Value = "'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D41").Value = "'0.9207"
$ws.Range("E41").Value = "'  +0.87%  "

$ws.Range("B42").Value = "'RenderToken"
$ws.Range("C42").Value = "'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D42").Value = "'1.981"
$ws.Range("E42").Value = "'  -4.15%  "

$ws.Range("B43").Value = "'Quant"
$ws.Range("C43").Value = "'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D43").Value = "'106.00"
$ws.Range("E43").Value = "'  +0.02%  "

$ws.Range("B44").Value = "'PaxDollar"
$ws.Range("C44").Value = "'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
$ws.Range("D44").Value = "'0.9998"
$ws.Range("E44").Value = "'  +0.14%  "

$ws.Range("B45").Value = "'TheSandbox"
$ws.Range("C45").Value = "'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D45").Value = "'0.4122"
$ws.Range("E45").Value = "'  -3.15%  "

$ws.Range("B46").Value = "'FraxShare"
$ws.Range("C46").Value = "'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D46").Value = "'5.603"
$ws.Range("E46").Value = "'  -2.64%  "

$ws.Range("B47").Value = "'Aptos"
$ws.Range("C47").Value = "'https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D47").Value = "'7.135"
$ws.Range("E47").Value = "'  -2.69%  "

$ws.Range("B48").Value = "'Aave"
$ws.Range("C48").Value = "'https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D48").Value = "'61.21"
$ws.Range("E48").Value = "'  -4.04%  "

$ws.Range("B49").Value = "'Algorand"
$ws.Range("C49").Value = "'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D49").Value = "'0.1212"
$ws.Range("E49").Value = "'  -5.40%  "

$ws.Range("B50").Value = "'EnergySwap"
$ws.Range("C50").Value = "'https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D50").Value = "'8.709"
$ws.Range("E50").Value = "'  -2.00%  "

$ws.Range("B51").Value = "'Elrond"
$ws.Range("C51").Value = "'https://coinranking.com/coin/omwkOTglq+elrond-egld"
$ws.Range("D51").Value = "'33.57"
$ws.Range("E51").Value = "'  -0.11%  "
